# refactorización features y data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

# Update test data for e-prepago data-driven sheet
$ws.Range("B3").Value = 48646663
$ws.Range("D3").Value = "autotest11"

# Update the selected/active cell to match refactored state
$ws.Range("A3").Select()
